# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# Column G header is "K"; only the cached numeric values change - each row's
# K value is recomputed/rewritten with the new source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    7  = 1
    8  = 3
    9  = 0
    10 = 0
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 1
    26 = 2
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 2
    32 = 2
    34 = 2
    35 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
